$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast error values (Bugfixed evaluation and simulated rt_data for components)
$data = @{
    2  = @(-0.1286665469352575, 2.584314119446431, 11.86673044746297, 3.444812106263994, 3.475993634475042, 52)
    3  = @(0.208472532147568, 1.905008350356748, 5.6589047777894, 2.378845261421894, 2.393272416744899, 51)
    4  = @(0.3820242035387909, 1.902010132220234, 5.867547615856401, 2.422302131414742, 2.416272454685605, 50)
    5  = @(0.191239345277837, 1.716821400939533, 5.499592484960565, 2.345120995803962, 2.36634644499125, 41)
    6  = @(0.2671928775352811, 1.741437505328719, 5.032127257509086, 2.243240347691055, 2.264087689223547, 31)
    7  = @(0.3126204252219001, 1.694986524670226, 4.73978801325031, 2.177105420793929, 2.192679688080215, 29)
    8  = @(0.3719975688557745, 1.867752137679124, 5.380359802164317, 2.319560260515841, 2.333150688324054, 27)
    9  = @(0.4130801636093303, 2.085153282719833, 6.310345893095409, 2.512040185406159, 2.545742663513931, 19)
    10 = @(1.987269086558114, 1.995519100223671, 6.088527956227622, 2.46749426670613, 1.527668288782175, 12)
    11 = @(1.658904638620646, 1.679380840373142, 3.778169217003525, 1.943751325916853, 1.132588085407946, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
